# Update cryptos list - GitHub Actions scheduled refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) holds numeric-looking values that must stay as
# literal text (as in the source data, e.g. "65.401.16" / "138.40").
# Force text formatting on every Price cell we touch so Excel's COM value
# setter doesn't silently coerce them into numbers (which would also drop
# meaningful trailing zeros).
$priceCells = @("D2","D3","D5","D6","D7","D10","D11","D13","D15","D16","D18","D21","D22","D24","D25","D27","D28","D29","D31","D32","D34","D37","D38","D39","D43","D44","D47","D50","D51")
foreach ($cell in $priceCells) {
    $ws.Range($cell).NumberFormat = "@"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "65.423.86"
$ws.Range("E2").Value = "  +0.31%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.543.77"
$ws.Range("E3").Value = "  +4.00%  "

# Row 5 - BNB
$ws.Range("D5").Value = "600.64"
$ws.Range("E5").Value = "  +3.07%  "

# Row 6 - Solana
$ws.Range("D6").Value = "138.40"
$ws.Range("E6").Value = "  +1.46%  "

# Row 7 - LidoStakedEther
$ws.Range("D7").Value = "3.542.23"
$ws.Range("E7").Value = "  +3.94%  "

# Row 9 - XRP
$ws.Range("E9").Value = "  -0.25%  "

# Row 10 - Dogecoin
$ws.Range("D10").Value = "0.124"
$ws.Range("E10").Value = "  +3.76%  "

# Row 11 - Toncoin
$ws.Range("D11").Value = "6.89"

# Row 12 - Cardano
$ws.Range("E12").Value = "  +3.67%  "

# Row 13 - WrappedliquidstakedEther2.0
$ws.Range("D13").Value = "4.146.08"
$ws.Range("E13").Value = "  +4.03%  "

# Row 14 - ShibaInu
$ws.Range("E14").Value = "  +3.75%  "

# Row 15 - Avalanche
$ws.Range("D15").Value = "27.22"
$ws.Range("E15").Value = "  +5.34%  "

# Row 16 - WrappedEther
$ws.Range("D16").Value = "3.546.64"
$ws.Range("E16").Value = "  +4.02%  "

# Row 17 - TRON
$ws.Range("E17").Value = "  +1.45%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "65.283.79"
$ws.Range("E18").Value = "  +0.05%  "

# Row 19 - Uniswap
$ws.Range("E19").Value = "  +5.64%  "

# Row 20 - Polkadot
$ws.Range("E20").Value = "  +1.94%  "

# Row 21 - Chainlink
$ws.Range("D21").Value = "14.27"
$ws.Range("E21").Value = "  +5.78%  "

# Row 22 - BitcoinCash
$ws.Range("D22").Value = "392.50"
$ws.Range("E22").Value = "  +2.63%  "

# Row 23 - Polygon
$ws.Range("E23").Value = "  +3.88%  "

# Row 24 - WrappedeETH
$ws.Range("D24").Value = "3.691.55"
$ws.Range("E24").Value = "  +4.15%  "

# Row 25 - Litecoin
$ws.Range("D25").Value = "73.65"
$ws.Range("E25").Value = "  +1.49%  "

# Row 26 - Dai
$ws.Range("E26").Value = "  +0.10%  "

# Row 27 - PEPE
$ws.Range("D27").Value = "0.0000115"
$ws.Range("E27").Value = "  +10.65%  "

# Row 28 - RenderToken
$ws.Range("D28").Value = "7.80"
$ws.Range("E28").Value = "  +11.27%  "

# Row 29 - Binance-PegBSC-USD
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  -0.12%  "

# Row 30 - PancakeSwap
$ws.Range("E30").Value = "  +3.98%  "

# Row 31 - InternetComputer(DFINITY)
$ws.Range("D31").Value = "8.17"
$ws.Range("E31").Value = "  +2.09%  "

# Row 32 - RenzoRestakedETH
$ws.Range("D32").Value = "3.556.59"
$ws.Range("E32").Value = "  +4.23%  "

# Row 34 - EthereumClassic
$ws.Range("D34").Value = "23.84"
$ws.Range("E34").Value = "  +4.85%  "

# Row 35 - Kaspa
$ws.Range("E35").Value = "  +2.33%  "

# Row 36 - Fetch.AI
$ws.Range("E36").Value = "  +15.66%  "

# Row 37 - was Aptos, now ImmutableX
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").Value = "1.57"
$ws.Range("E37").Value = "  +8.71%  "

# Row 38 - Monero
$ws.Range("D38").Value = "169.89"
$ws.Range("E38").Value = "  -0.56%  "

# Row 39 - was ImmutableX, now Aptos
$ws.Range("B39").Value = "Aptos"
$ws.Range("C39").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D39").Value = "6.95"
$ws.Range("E39").Value = "  +3.35%  "

# Row 40 - NEARProtocol
$ws.Range("E40").Value = "  +6.46%  "

# Row 41 - Hedera
$ws.Range("E41").Value = "  +6.55%  "

# Row 42 - Mantle
$ws.Range("E42").Value = "  +1.68%  "

# Row 43 - EnergySwap
$ws.Range("D43").Value = "26.56"
$ws.Range("E43").Value = "  +20.31%  "

# Row 44 - OKB
$ws.Range("D44").Value = "42.67"
$ws.Range("E44").Value = "  -1.79%  "

# Row 45 - FirstDigitalUSD
$ws.Range("E45").Value = "  -0.08%  "

# Row 46 - Filecoin
$ws.Range("E46").Value = "  +1.60%  "

# Row 47 - Stacks
$ws.Range("D47").Value = "1.69"
$ws.Range("E47").Value = "  +6.54%  "

# Row 48 - ONDO
$ws.Range("E48").Value = "  +10.89%  "

# Row 49 - Cosmos
$ws.Range("E49").Value = "  +5.10%  "

# Row 50 - Maker
$ws.Range("D50").Value = "2.396.70"
$ws.Range("E50").Value = "  +10.39%  "

# Row 51 - was Bittensor, now LidoDAOToken
$ws.Range("B51").Value = "LidoDAOToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D51").Value = "2.38"
$ws.Range("E51").Value = "  +20.98%  "
